$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (51 -> 52 chars). Runtime applies an MDW-based offset of +5/6
# to ColumnWidth input before storing <col width>, so feed 52 - 5/6 to land on exactly 52.
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664

# Drop all existing hyperlinks first -- rows are being reshuffled/renumbered, and
# leftover relationships at the old F-cell refs would collide with the new ones
# added below for the same cells.
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2026-01-23 12:42:09'
$ws.Range("B2").Value = '製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5460562'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5460562')
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 435
$ws.Range("H2").Value = '🔥AI,Ai ◆ツール,開発'

# Row 3
$ws.Range("A3").Value = '2026-01-23 12:42:09'
$ws.Range("B3").Value = '【急募】AIシステム開発で情報抽出・転記のプロを求む!'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5477580'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5477580')
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value = 403
$ws.Range("H3").Value = '🔥AI,Ai ◆開発,システム開発'

# Row 4
$ws.Range("A4").Value = '2026-01-23 12:42:09'
$ws.Range("B4").Value = '【フルスタックエンジニア募集】AWS構築+Pythonバックエンド開発'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5475657'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5475657')
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("G4").Value = 260
$ws.Range("H4").Value = '🔥Python ◆開発'

# Row 5
$ws.Range("A5").Value = '2026-01-23 12:42:09'
$ws.Range("B5").Value = '※急募:Next.jsによる業務アプリの開発(+Flutter)'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5477335'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5477335')
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("G5").Value = 225
$ws.Range("H5").Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 6
$ws.Range("A6").Value = '2026-01-23 12:42:09'
$ws.Range("B6").Value = '施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5460563'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5460563')
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value = 220
$ws.Range("H6").Value = '◆開発,システム開発 ◇管理'

# Row 7
$ws.Range("A7").Value = '2026-01-23 12:42:09'
$ws.Range("B7").Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5477338'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5477338')
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 218
$ws.Range("H7").Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 8
$ws.Range("A8").Value = '2026-01-23 12:42:09'
$ws.Range("B8").Value = '【急募】データ管理ツールの開発をお手伝いください!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5477312'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5477312')
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("G8").Value = 163
$ws.Range("H8").Value = '◆ツール,開発 ◇管理'

# Row 9
$ws.Range("A9").Value = '2026-01-23 12:42:09'
$ws.Range("B9").Value = '【急募】魅力的なWebシステム開発のパートナーを探しています!'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5477481'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5477481')
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("G9").Value = 118
$ws.Range("H9").Value = '◆開発,システム開発'

# Row 10
$ws.Range("A10").Value = '2026-01-23 12:42:09'
$ws.Range("B10").Value = '自動化システム'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5477084'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5477084')
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("G10").Value = 110
$ws.Range("H10").Value = '◆自動化'

# Row 11
$ws.Range("A11").Value = '2026-01-23 12:42:09'
$ws.Range("B11").Value = '【急募】宿泊業向けクチコミ対策SaaSのMVP開発'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5466852'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5466852')
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("G11").Value = 75
$ws.Range("H11").Value = '◆開発'

# Row 12
$ws.Range("A12").Value = '2026-01-23 12:42:09'
$ws.Range("B12").Value = 'Keepaの通知からAmazonで自動購入するシステムの開発依頼の仕事'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5477013'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5477013')
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("G12").Value = 75
$ws.Range("H12").Value = '◆開発'

# Row 13
$ws.Range("A13").Value = '2026-01-23 12:42:09'
$ws.Range("B13").Value = 'Keepaの通知からAmazonで自動購入するシステムの開発依頼の仕事'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5477036'
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5477036')
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("G13").Value = 75
$ws.Range("H13").Value = '◆開発'

# Row 14
$ws.Range("A14").Value = '2026-01-23 12:42:09'
$ws.Range("B14").Value = '【急募】見積依頼集約と遅延防止のMicrosoft365システム構築'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5477550'
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5477550')
$ws.Range("F14").Style = "Hyperlink"
$ws.Range("G14").Value = 33

# Row 15
$ws.Range("A15").Value = '2026-01-23 12:42:09'
$ws.Range("B15").Value = '以前1/60秒単位のリピートタイマーを作成。これを2連(ダブルカウントダウンタイマー)から6連に。'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5477366'
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5477366')
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("G15").Value = 18
